# Auto-generated PowerShell COM-interop script
# Implements: insert a new "2022-Q3" worksheet (with fund holding detail data)
# right after "总计" and before "2022-Q2"; update the "总计" summary sheet with
# a new leading row for 2022-Q3 (shifting the existing quarters down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (inline/shared string), even when
# it looks like a number (e.g. fund codes such as "010490", or decimal
# figures like "41.17"). Plain `.Value = "41.17"` auto-coerces to a number
# (losing leading zeros / matching the diff's inlineStr cells), so we briefly
# format the cell as Text, assign the value, then clear the format again so
# no residual style/quotePrefix is left behind (matches the target XML,
# which has no `s` attribute and no `quotePrefix` on these cells).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# Helper: stamp the shared "index / bold+border" style used on column A and
# on header rows (style index 2 in the original styles.xml) onto a cell by
# copying it from a cell that already carries that exact style, then
# overwrite the copied value. Using Copy() reuses the existing style slot
# instead of synthesizing a near-duplicate one.
# ---------------------------------------------------------------------------
function Set-StyledValue {
    param($StyleSource, $Cell, $Value)
    $StyleSource.Copy($Cell)
    $Cell.Value = $Value
}

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right after "总计" (so order becomes
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Style-reference cells (already carrying the bold/border style used
# throughout these sheets) that we will clone styling from without
# creating new style entries.
$styleRefHeader = $wb.Worksheets.Item(3).Range("B1")   # "2022-Q2" header cell, style 2
$styleRefIndex  = $wb.Worksheets.Item(3).Range("A2")   # "2022-Q2" index cell, style 2

# ---------------------------------------------------------------------------
# 2) Populate header row 1 (B1:H1) of the new sheet.
# ---------------------------------------------------------------------------
Set-StyledValue $styleRefHeader $q3Sheet.Range("B1") '基金代码'
$q3Sheet.Range("B1").NumberFormat = "@"
Set-StyledValue $styleRefHeader $q3Sheet.Range("C1") '基金名称'
$q3Sheet.Range("C1").NumberFormat = "@"
Set-StyledValue $styleRefHeader $q3Sheet.Range("D1") '基金规模'
$q3Sheet.Range("D1").NumberFormat = "@"
Set-StyledValue $styleRefHeader $q3Sheet.Range("E1") '股票总仓位'
$q3Sheet.Range("E1").NumberFormat = "@"
Set-StyledValue $styleRefHeader $q3Sheet.Range("F1") '仓位占比'
$q3Sheet.Range("F1").NumberFormat = "@"
Set-StyledValue $styleRefHeader $q3Sheet.Range("G1") '持有市值(亿元)'
$q3Sheet.Range("G1").NumberFormat = "@"
Set-StyledValue $styleRefHeader $q3Sheet.Range("H1") '仓位排名'
$q3Sheet.Range("H1").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3) Populate the 43 data rows (rows 2-44).
# ---------------------------------------------------------------------------
$q3Rows = @(
    @(0, '206009', '鹏华新兴产业混合', '41.17', '89.21', '3.72', '1.5315', 10),
    @(1, '398021', '中海能源策略混合', '21.73', '88.03', '5.90', '1.2821', 2),
    @(2, '010490', '鹏华高质量增长混合A', '12.74', '93.98', '9.52', '1.2128', 1),
    @(3, '206002', '鹏华精选成长混合A', '12.50', '93.83', '5.41', '0.6762', 4),
    @(4, '398051', '中海环保新能源混合', '22.19', '67.74', '2.96', '0.6568', 10),
    @(5, '009984', '鹏华启航混合', '12.45', '84.87', '4.19', '0.5217', 7),
    @(6, '011956', '鹏华新能源精选混合A', '7.18', '85.90', '6.08', '0.4365', 2),
    @(7, '006976', '鹏华核心优势混合', '4.55', '89.77', '8.63', '0.3927', 1),
    @(8, '009023', '鹏华稳健回报混合', '4.12', '94.39', '8.39', '0.3457', 2),
    @(9, '011460', '鹏华创新成长混合A', '10.08', '87.81', '3.08', '0.3105', 10),
    @(10, '011957', '鹏华新能源精选混合C', '4.87', '85.90', '6.08', '0.2961', 2),
    @(11, '398061', '中海消费混合', '3.91', '85.30', '6.23', '0.2436', 2),
    @(12, '370024', '上投摩根核心优选混合A', '8.10', '78.00', '2.01', '0.1628', 10),
    @(13, '009990', '华泰柏瑞品质优选混合A', '9.06', '68.27', '1.76', '0.1595', 10),
    @(14, '008811', '鹏华科技创新混合', '3.02', '88.94', '4.81', '0.1453', 3),
    @(15, '000431', '鹏华品牌传承混合', '4.52', '75.03', '2.58', '0.1166', 9),
    @(16, '004986', '鹏华策略回报灵活配置混合', '3.12', '87.35', '3.19', '0.0995', 4),
    @(17, '952035', '国泰君安君得诚混合', '2.35', '80.66', '4.18', '0.0982', 5),
    @(18, '206012', '鹏华价值精选股票', '2.64', '86.31', '3.23', '0.0853', 5),
    @(19, '000166', '中海信息产业精选混合', '0.77', '89.31', '7.10', '0.0547', 1),
    @(20, '009991', '华泰柏瑞品质优选混合C', '2.49', '68.27', '1.76', '0.0438', 10),
    @(21, '010491', '鹏华高质量增长混合C', '0.44', '93.98', '9.52', '0.0419', 1),
    @(22, '006526', '鹏华优选回报灵活配置混合A', '0.94', '80.05', '2.82', '0.0265', 10),
    @(23, '016562', '鹏华精选成长混合C', '0.34', '93.83', '5.41', '0.0184', 4),
    @(24, '000354', '长盛城镇化主题混合', '0.34', '92.91', '5.03', '0.0171', 8),
    @(25, '004258', '国寿安保稳嘉混合A', '2.32', '20.12', '0.69', '0.0160', 10),
    @(26, '000743', '红塔红土盛世普益灵活配置混合', '1.12', '20.22', '1.29', '0.0144', 4),
    @(27, '011461', '鹏华创新成长混合C', '0.37', '87.81', '3.08', '0.0114', 10),
    @(28, '004301', '国寿安保稳信混合A', '1.50', '22.10', '0.71', '0.0106', 9),
    @(29, '620004', '金元顺安价值增长混合', '0.34', '74.99', '1.99', '0.0068', 5),
    @(30, '008093', '同泰慧选混合A', '0.21', '63.44', '2.85', '0.0060', 7),
    @(31, '008094', '同泰慧选混合C', '0.13', '63.44', '2.85', '0.0037', 7),
    @(32, '009027', '浦银安盛安远回报一年持有期混合A', '0.57', '21.48', '0.53', '0.0030', 10),
    @(33, '002023', '红塔红土稳健回报灵活配置混合A', '0.09', '62.51', '3.21', '0.0029', 8),
    @(34, '004276', '浦银安盛安和回报定期开放混合A', '0.39', '20.97', '0.53', '0.0021', 9),
    @(35, '012997', '鹏华优选回报灵活配置混合C', '0.04', '80.05', '2.82', '0.0011', 10),
    @(36, '009028', '浦银安盛安远回报一年持有期混合C', '0.10', '21.48', '0.53', '0.0005', 10),
    @(37, '004277', '浦银安盛安和回报定期开放混合C', '0.03', '20.97', '0.53', '0.0002', 9),
    @(38, '004302', '国寿安保稳信混合C', '0.02', '22.10', '0.71', '0.0001', 9),
    @(39, '015406', '国寿安保稳信混合E', '0.01', '22.10', '0.71', '0.0001', 9),
    @(40, '004259', '国寿安保稳嘉混合C', '0.00', '20.12', '0.69', 0, 10),
    @(41, '015057', '上投摩根核心优选混合C', '0.00', '78.00', '2.01', 0, 10),
    @(42, '002024', '红塔红土稳健回报灵活配置混合C', '0.00', '62.51', '3.21', 0, 8)
)

foreach ($r in $q3Rows) {
    $rowNum = [int]$r[0] + 2
    Set-StyledValue $styleRefIndex $q3Sheet.Cells.Item($rowNum, 1) ([int]$r[0])
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 2) $r[1]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 3) $r[2]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 4) $r[3]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 5) $r[4]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 6) $r[5]
    $gValue = $r[6]
    if ($gValue -is [int]) {
        $q3Sheet.Cells.Item($rowNum, 7).Value = $gValue
    } else {
        Set-TextValue $q3Sheet.Cells.Item($rowNum, 7) $gValue
    }
    $q3Sheet.Cells.Item($rowNum, 8).Value = $r[7]
}

# ---------------------------------------------------------------------------
# 4) Update the "总计" summary sheet: insert a new leading data row for
#    2022-Q3 (index 0) and shift the existing quarters' index down by one.
# ---------------------------------------------------------------------------
$totalRows = @(
    @(0, '2022-Q3', 43, 9.050000000000001),
    @(1, '2022-Q2', 14, 4.22),
    @(2, '2022-Q1', 6, 1.52),
    @(3, '2021-Q4', 13, 1.49),
    @(4, '2021-Q3', 2, 0.03)
)

$totalStyleRefIndex = $totalSheet.Range("A2")
foreach ($r in $totalRows) {
    $rowNum = [int]$r[0] + 2
    Set-StyledValue $totalStyleRefIndex $totalSheet.Cells.Item($rowNum, 1) ([int]$r[0])
    Set-TextValue $totalSheet.Cells.Item($rowNum, 2) $r[1]
    $totalSheet.Cells.Item($rowNum, 3).Value = $r[2]
    $totalSheet.Cells.Item($rowNum, 4).Value = $r[3]
}

Write-Output "done"
